# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Illes Balears*" (row 26) and "Illes Balears" (row 27)
$ws.Range("A26").Value = "Illes Balears"
$ws.Range("A27").Value = "Illes Balears*"

# Swap "Huelva" (row 52) and "Melilla" (row 53)
$ws.Range("A52").Value = "Melilla"
$ws.Range("A53").Value = "Huelva"

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 02:46"
